$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 149, pushing the existing rows 149:282 down to 150:283.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with a fresh "Apio" price entry for
# Macroferia Regional de Talca / Maule, duplicating the metadata that was
# previously on row 149 (now row 150) but with a new reporting date.
$ws.Range("A149").Value = 5
$ws.Range("B149").Value = "Macroferia Regional de Talca"
$ws.Range("C149").Value = "Maule"
$ws.Range("D149").Value = 45040
$ws.Range("E149").Value = 7
$ws.Range("F149").Value = 100112017
$ws.Range("G149").Value = "Apio"
$ws.Range("H149").Value = "Americana (o)"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 600
$ws.Range("K149").Value = 7000
$ws.Range("L149").Value = 7000
$ws.Range("M149").Value = 7000
$ws.Range("N149").Value = "$/docena de matas"
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 1167
$ws.Range("Q149").Value = 6
$ws.Range("R149").Value = "Hortaliza"
